$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number that was bumped by one day
# (2025-06-03 -> 2025-06-04) for every data row (rows 2 through 43).
for ($row = 2; $row -le 43; $row++) {
    $ws.Cells.Item($row, 3).Value = 45812
}
